# Hierarchy by children is working
$wb = $excel.ActiveWorkbook

# Duplicate the existing sheet so the new "Per Child" sheet starts out with
# the same formatting (page margins, row height, etc.) as "Per Parent"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)

$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "Per Parent"
$ws2.Name = "Per Child"

# Header row
$ws2.Range("C1").Value = "Child"

# Data rows - same Id/Name as "Per Parent"; the third column now lists the
# comma-separated ids of the children instead of a single parent id.
$ws2.Range("C2").Value = "2,6,7,8,10"
$ws2.Range("C3").Value = 3
$ws2.Range("C4").Value = 4
$ws2.Range("C5").Value = 5
$ws2.Range("C6").Value = ""
$ws2.Range("C7").Value = ""
$ws2.Range("C8").Value = ""
$ws2.Range("C9").Value = ""
$ws2.Range("C10").Value = ""
$ws2.Range("C11").Value = 9

# Restore selections on both sheets: full table selected on "Per Parent",
# and an arbitrary cell below the table selected (and active) on "Per Child"
[void]$ws1.Range("A1:C11").Select()
[void]$ws2.Range("J12").Select()
